$p = $ppt.ActivePresentation

# Slide 1
$s1 = $p.Slides.Item(1)

# TextBox 5 ("Executive Gold" title) - sz 4800 b1 64503C
$tr1 = $s1.Shapes.Item("TextBox 5").TextFrame.TextRange
$tr1.ParagraphFormat.Alignment = 1  # ppAlignLeft
$tr1.Font.Name = "Arial"

# TextBox 6 ("Premium & Luxury" subtitle) - sz 2000 b0 8B7555
$tr2 = $s1.Shapes.Item("TextBox 6").TextFrame.TextRange
$tr2.ParagraphFormat.Alignment = 1  # ppAlignLeft
$tr2.Font.Name = "Arial"

# Slide 2
$s2 = $p.Slides.Item(2)

# TextBox 3 ("Premium Features" title) - sz 3200 b1 FFFFFF
$tr3 = $s2.Shapes.Item("TextBox 3").TextFrame.TextRange
$tr3.ParagraphFormat.Alignment = 1  # ppAlignLeft
$tr3.Font.Name = "Arial"

# TextBox 5 (bullet list) - sz 2000 b0 50463C
$tr4 = $s2.Shapes.Item("TextBox 5").TextFrame.TextRange
$tr4.ParagraphFormat.Alignment = 1  # ppAlignLeft
$tr4.Font.Name = "Arial"
